$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: new entry (No=8, Date=2022-01-12, Application=RPA RLOGIC,
# Task text, % of completion=1 (100%), Status=Completed).
# Enter the new Task text first so it becomes the first newly-added
# shared string.
$ws.Range("D15").Value = "1. Correction  received for  the Accounting Statements for the three centers and completed the correction and shared to Rahaman san to verify"
$ws.Range("A15").Value = 8

# Copy the date/percent number formats from existing rows so the same
# (reused) style indexes end up being referenced, instead of creating
# brand new number formats.
$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B15").Value = 44573

$ws.Range("C15").Value = "RPA RLOGIC"

$ws.Range("E13").Copy()
$ws.Range("E15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E15").Value = 1

$ws.Range("F15").Value = "Completed"

# Row 14: set C14 to "GSS" (second newly-added shared string)
$ws.Range("C14").Value = "GSS"

# Row 16: continuation entry (Application=GSS, Task text,
# % of completion=0.5 (50%), Status=WIP)
$ws.Range("C16").Value = "GSS"
$ws.Range("D16").Value = "2. Uploading the pdf files task is work in progress for ESA  and it is around 50% has been completed"

$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E16").Value = 0.5

$ws.Range("F16").Value = "WIP"

# Update the active selection to D18, matching the saved view state
$ws.Range("D18").Select()
